# Swap the order of the two comma-separated names/emails in the
# "Recorded By" (column G) cells that currently read
# "dnasr281@gmail.com, System" or "dnasr281@gmail.com, admin@admin.com"
# so that they become "System, dnasr281@gmail.com" /
# "admin@admin.com, dnasr281@gmail.com" respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G - "Recorded By"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($val -ne $null) {
        if ($val -eq "dnasr281@gmail.com, System" -or $val -eq "dnasr281@gmail.com, admin@admin.com") {
            $parts = $val -split ", ", 2
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value = $newVal
        }
    }
}
